# ============================================================================
# This script turns the single-sheet workbook ("ODI Batting") into a
# three-sheet workbook:
#   1. "Player Info"        (new)
#   2. "ODI Batting"        (existing sheet, D column changed from a full
#                             URL to just the bare MatchCode, header renamed
#                             from MATCH_CARD_LINK to MATCH_CODE)
#   3. "ODI Batting Extra"  (new)
# ============================================================================

$wb = $excel.ActiveWorkbook
$odiBatting = $wb.Worksheets.Item(1)

# ----------------------------------------------------------------------
# Helper: make the header row of a range bold / centered / top-aligned /
# thin-bordered, matching the look of the existing "ODI Batting" header.
# ----------------------------------------------------------------------
function Format-HeaderRange($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# ----------------------------------------------------------------------
# 1) Add the "Player Info" sheet before the existing "ODI Batting" sheet
# ----------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $playerInfoHeaders.Length; $col++) {
    $playerInfo.Cells.Item(1, $col).Value = $playerInfoHeaders[$col - 1]
}
Format-HeaderRange ($playerInfo.Range("A1:D1"))

$playerInfo.Cells.Item(2, 1).Value = "'5665"
$playerInfo.Cells.Item(2, 2).Value = "Brandon Alexander King"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# ----------------------------------------------------------------------
# 2) Add the "ODI Batting Extra" sheet after the existing "ODI Batting"
#    sheet
#    NOTE: worksheet references returned earlier become stale once the
#    sheet collection is mutated (this engine resolves object handles by
#    position, not stable identity), so re-fetch "ODI Batting" by name.
# ----------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$extra = $wb.Worksheets.Add($null, $odiBatting)
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $extraHeaders.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $extraHeaders[$col - 1]
}
Format-HeaderRange ($extra.Range("A1:F1"))

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    ,@("4394", "4", "0", "0", $null, "NO")
    ,@("4397", "4", "4", "1", "19.10%", "NO")
    ,@("4533", $null, $null, $null, $null, "NO")
    ,@("4535", "2", "2", "1", "9.33%", "NO")
    ,@("4536", "2", "2", "0", "8.28%", "NO")
    ,@("4577", "5", "5", "2", "23.29%", "NO")
    ,@("4580", "5", "9", "3", "41.94%", "YES")
    ,@("4583", "5", "1", "0", "3.25%", "NO")
    ,@("4586", $null, $null, $null, $null, "NO")
    ,@("4590", $null, $null, $null, $null, "NO")
    ,@("4606", $null, $null, $null, $null, "NO")
    ,@("4611", "4", "0", "0", "10.19%", "NO")
    ,@("4616", "2", "1", "0", "4.49%", "NO")
    ,@("4621", "4", "2", "2", "17.70%", "NO")
    ,@("4623", $null, $null, $null, $null, "NO")
    ,@("4624", "4", "5", "1", "30.66%", "NO")
    ,@("4639", "4", "0", "0", "1.24%", "NO")
    ,@("4642", $null, $null, $null, $null, "NO")
    ,@("4727", $null, $null, $null, $null, "NO")
    ,@("4731", "1", "11", "1", "27.69%", "NO")
)

$rowIdx = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($rowIdx, 1).Value = "'" + $row[0]
    if ($row[1] -ne $null) {
        $extra.Cells.Item($rowIdx, 2).Value = [int]$row[1]
    }
    if ($row[2] -ne $null) {
        $extra.Cells.Item($rowIdx, 3).Value = "'" + $row[2]
    }
    if ($row[3] -ne $null) {
        $extra.Cells.Item($rowIdx, 4).Value = "'" + $row[3]
    }
    if ($row[4] -ne $null) {
        $extra.Cells.Item($rowIdx, 5).Value = "'" + $row[4]
    }
    $extra.Cells.Item($rowIdx, 6).Value = $row[5]
    $rowIdx = $rowIdx + 1
}

# ----------------------------------------------------------------------
# 3) Update the existing "ODI Batting" sheet: rename MATCH_CARD_LINK ->
#    MATCH_CODE and replace the full scorecard URL with the bare
#    MatchCode number in column D.
# ----------------------------------------------------------------------
$odiBatting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @("4379", "4391", "4394", "4397", "4533", "4535", "4536", "4577", "4580", "4583", "4586", "4590", "4606", "4611", "4616", "4621", "4623", "4624", "4639", "4642", "4727", "4731")

$rowIdx = 2
foreach ($code in $matchCodes) {
    $odiBatting.Cells.Item($rowIdx, 4).Value = "'" + $code
    $rowIdx = $rowIdx + 1
}
